# Auto-generated edit script applying the Titan_Profits.xlsx diff
# (sheet names in this workbook: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# ALC row 129 (diff hunk @ 7079)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1120.6875
$ws.Range("I129").Value = 344.75
$ws.Range("K129").Value = 1034.25
$ws.Range("M129").Value = 3965.75

# ALC row 132 (diff hunk @ 7232)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 40812.383
$ws.Range("I132").Value = 42364.88
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 127094.64
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -124564.64
$ws.Range("N132").Value = -11060

# ALC row 137 (diff hunk @ 7480)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 23257000
$ws.Range("I137").Value = 29412614
$ws.Range("K137").Value = 88237842
$ws.Range("M137").Value = -88235292

# ARM row 2 (diff hunk @ 7828)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35411.55
$ws.Range("I2").Value = 38004.258
$ws.Range("K2").Value = 38004.258
$ws.Range("M2").Value = -37891.258

# ARM row 32 (diff hunk @ 9322)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16312.161
$ws.Range("I32").Value = 3467.8076
$ws.Range("J32").Value = 83102.8
$ws.Range("K32").Value = 3467.8076
$ws.Range("L32").Value = 83102.8
$ws.Range("M32").Value = -3180.8076
$ws.Range("N32").Value = -83676.8

# ARM row 45 (diff hunk @ 9959)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 976.6667
$ws.Range("I45").Value = 981.6667
$ws.Range("K45").Value = 981.6667
$ws.Range("M45").Value = -604.6667

# ARM row 61 (diff hunk @ 10740)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1789.55
$ws.Range("I61").Value = 1288.0588
$ws.Range("J61").Value = 4631.3335
$ws.Range("K61").Value = 1288.0588
$ws.Range("L61").Value = 4631.3335
$ws.Range("M61").Value = -1076.0588
$ws.Range("N61").Value = -5055.3335

# ARM row 74 (diff hunk @ 11374)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4702.9443
$ws.Range("I74").Value = 850.04346
$ws.Range("J74").Value = 11519.615
$ws.Range("K74").Value = 850.04346
$ws.Range("L74").Value = 11519.615
$ws.Range("M74").Value = 23.95654000000002
$ws.Range("N74").Value = -13267.615

# ARM row 77 (diff hunk @ 11518)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4702.9443
$ws.Range("I77").Value = 850.04346
$ws.Range("J77").Value = 11519.615
$ws.Range("K77").Value = 4250.2173
$ws.Range("L77").Value = 57598.075
$ws.Range("M77").Value = 117.7826999999997
$ws.Range("N77").Value = -66334.075

# ARM row 116 (diff hunk @ 13438)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 35411.55
$ws.Range("I116").Value = 38004.258
$ws.Range("K116").Value = 38004.258
$ws.Range("M116").Value = -35710.258

# ARM row 122 (diff hunk @ 13732)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6967.636
$ws.Range("I122").Value = 6390.6665
$ws.Range("J122").Value = 7660
$ws.Range("K122").Value = 19171.9995
$ws.Range("L122").Value = 22980
$ws.Range("M122").Value = -16721.9995
$ws.Range("N122").Value = -27880

# ARM row 132 (diff hunk @ 14219)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2781.3447
$ws.Range("I132").Value = 2714.36
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 8143.08
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -5613.08
$ws.Range("N132").Value = -14660

# ARM row 136 (diff hunk @ 14418)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1789.55
$ws.Range("I136").Value = 1288.0588
$ws.Range("J136").Value = 4631.3335
$ws.Range("K136").Value = 3864.1764
$ws.Range("L136").Value = 13894.0005
$ws.Range("M136").Value = -1314.1764
$ws.Range("N136").Value = -18994.0005

# BSM row 3 (diff hunk @ 14855)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35411.55
$ws.Range("I3").Value = 38004.258
$ws.Range("K3").Value = 38004.258
$ws.Range("M3").Value = -37890.258

# BSM row 20 (diff hunk @ 15709)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2049.9
$ws.Range("I20").Value = 2166.6667
$ws.Range("K20").Value = 2166.6667
$ws.Range("M20").Value = -1919.6667

# BSM row 134 (diff hunk @ 21247)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5552
$ws.Range("I134").Value = 5552
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 16656
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -14121
$ws.Range("N134").ClearContents()

# CRP row 28 (diff hunk @ 23022)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 30000
$ws.Range("J28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("N28").Value = -30490

# CRP row 58 (diff hunk @ 24501)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3106.5
$ws.Range("I58").Value = 1511.25
$ws.Range("J58").Value = 4701.75
$ws.Range("K58").Value = 1511.25
$ws.Range("L58").Value = 4701.75
$ws.Range("M58").Value = -1308.25
$ws.Range("N58").Value = -5107.75

# CRP row 136 (diff hunk @ 28323)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3106.5
$ws.Range("I136").Value = 1511.25
$ws.Range("J136").Value = 4701.75
$ws.Range("K136").Value = 4533.75
$ws.Range("L136").Value = 14105.25
$ws.Range("M136").Value = -1983.75
$ws.Range("N136").Value = -19205.25

# CUL row 17 (diff hunk @ 29467)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 75.5
$ws.Range("I17").Value = 61
$ws.Range("K17").Value = 183
$ws.Range("M17").Value = -14

# CUL row 25 (diff hunk @ 29871)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 540
$ws.Range("I25").Value = 540
$ws.Range("K25").Value = 1620
$ws.Range("M25").Value = -1451

# CUL row 30 (diff hunk @ 30128)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 540
$ws.Range("I30").Value = 540
$ws.Range("K30").Value = 1620
$ws.Range("M30").Value = -1518

# CUL row 34 (diff hunk @ 30327)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2607.2144
$ws.Range("J34").Value = 2461.4614
$ws.Range("L34").Value = 7384.3842
$ws.Range("N34").Value = -7552.3842

# CUL row 131 (diff hunk @ 35242)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6668157
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 7577382.5
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 22732147.5
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -22742227.5

# GSM row 122 (diff hunk @ 41743)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1926.8823
$ws.Range("I122").Value = 2213.0908
$ws.Range("J122").Value = 1402.1666
$ws.Range("K122").Value = 6639.2724
$ws.Range("L122").Value = 4206.4998
$ws.Range("M122").Value = -4189.2724
$ws.Range("N122").Value = -9106.4998

# LTW row 16 (diff hunk @ 43509)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8335092
$ws.Range("I16").Value = 10000911
$ws.Range("J16").Value = 5996
$ws.Range("K16").Value = 10000911
$ws.Range("L16").Value = 5996
$ws.Range("M16").Value = -10000741
$ws.Range("N16").Value = -6336

# LTW row 68 (diff hunk @ 46084)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2058
$ws.Range("I68").Value = 1995
$ws.Range("J68").Value = 2100
$ws.Range("K68").Value = 1995
$ws.Range("L68").Value = 2100
$ws.Range("M68").Value = -1246

# LTW row 71 (diff hunk @ 46231)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2058
$ws.Range("I71").Value = 1995
$ws.Range("J71").Value = 2100
$ws.Range("K71").Value = 9975
$ws.Range("L71").Value = 10500
$ws.Range("M71").Value = -6231
$ws.Range("N71").Value = -17988

# WVR row 98 (diff hunk @ 54535)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990

# WVR row 132 (diff hunk @ 56195)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7958.5835
$ws.Range("I132").Value = 8590.4
$ws.Range("J132").Value = 4799.5
$ws.Range("K132").Value = 25771.2
$ws.Range("L132").Value = 14398.5
$ws.Range("M132").Value = -23241.2
$ws.Range("N132").Value = -19458.5
